$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "LeaveEmpty" column to the end of the existing menu table (Table3).
$tbl = $ws.ListObjects.Item("Table3")
$tbl.ListColumns.Add() | Out-Null
$ws.Range("G1").Value = "LeaveEmpty"

# Fill in the missing Allergens value for the "Chicken Salad" row (row 4),
# matching the wrap-text style already used in row 3 for the same column.
$ws.Range("C4").Value = "No known priority allergens"
$ws.Range("C4").WrapText = $true

# Update the active selection to reflect the last edited cell.
$ws.Range("C4").Select()
